$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.88%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'31.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.55%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.069"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.52%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08120"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'9.59%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'16.67%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.812"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.56%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.826"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.63%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.05%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1761"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.89%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07475"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.97%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.31%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03032"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.29%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.61%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.10%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005882"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.41%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.555"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'2.80%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.250"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.03%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'-0.56%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1318"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.37%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'3.963"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-14.89%"
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'3.90%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04601"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.28%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001241"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.23%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004455"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.69%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-8.09%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0003408"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'82.08%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").Value = "'0.01767"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.76%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04526"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.10%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006864"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.26%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.38%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.90%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.009848"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-10.22%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006473"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.87%"
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'-0.29%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.008737"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'24.87%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'-57.44%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002094"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.29%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0001994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("E50").Style = "Normal"
